$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$CHECK = [char]0x2713

# ---------------------------------------------------------------------------
# 1. New column M: "Refined Tuning - Collect Data" mission mode
#    Header uses the same formatting as the other header cells (copy from L1).
# ---------------------------------------------------------------------------
$ws.Range("L1").Copy()
$ws.Range("M1").PasteSpecial(-4122)
$ws.Range("M1").Value = "Refined Tuning " + [char]0x2013 + " Collect Data"

# Rows that get a checkmark in the new column (style copied from a donor
# cell that already carries the "check" formatting, C10).
$checkRows = @(2,4,5,6,7,10,11,12,13,14,15,18)
foreach ($r in $checkRows) {
    $ws.Range("C10").Copy()
    $target = $ws.Cells.Item($r, 13)
    $target.PasteSpecial(-4122)
    $target.Value = $CHECK
}

# Rows that stay blank in the new column, but with plain formatting
# (style copied from donor cell B1).
$plainBlankRows = @(3,8,9,19,20,23,24,25,28,29,30,31,32,33)
foreach ($r in $plainBlankRows) {
    $ws.Range("B1").Copy()
    $target = $ws.Cells.Item($r, 13)
    $target.PasteSpecial(-4122)
    $target.ClearContents()
}

# Rows that are blank separator rows, with the separator-row formatting
# (style copied from donor cell C16).
$separatorBlankRows = @(16,17,21,22,26,27)
foreach ($r in $separatorBlankRows) {
    $ws.Range("C16").Copy()
    $target = $ws.Cells.Item($r, 13)
    $target.PasteSpecial(-4122)
    $target.ClearContents()
}

# ---------------------------------------------------------------------------
# 2. D10 now also shows the checkmark (new tuning availability for the
#    "RC Stabilize - Hover Thrust ID" mode).
# ---------------------------------------------------------------------------
$ws.Range("C10").Copy()
$ws.Range("D10").PasteSpecial(-4122)
$ws.Range("D10").Value = $CHECK

# ---------------------------------------------------------------------------
# 3. Row 26 clean-up: the stray checkmark in L26 is removed and the row
#    collapses back to the default (non-custom) row height, like its
#    neighbouring separator rows.
# ---------------------------------------------------------------------------
$ws.Range("K26").Copy()
$ws.Range("L26").PasteSpecial(-4122)
$ws.Range("L26").ClearContents()
$ws.Rows.Item(26).AutoFit()

# ---------------------------------------------------------------------------
# 4. Row height tweaks for rows 19 & 20.
# ---------------------------------------------------------------------------
$ws.Rows.Item(19).RowHeight = 53.45
$ws.Rows.Item(20).RowHeight = 53.45

# ---------------------------------------------------------------------------
# 5. Column width tweaks (G shrinks marginally, M widens for new content).
# ---------------------------------------------------------------------------
$ws.Columns.Item(7).ColumnWidth = 19.506666666666667
$ws.Columns.Item(13).ColumnWidth = 39.406666666666666

# ---------------------------------------------------------------------------
# 6. Sheet view: scroll position, zoom level and active selection.
# ---------------------------------------------------------------------------
$win = $excel.ActiveWindow
$win.DisplayGridlines = $true
$win.Zoom = 85
$win.ScrollColumn = 4
$win.ScrollRow = 1
$ws.Range("N34").Select()
